# Generate Report for Handoff
# Updates the localization-status report: flips the "Status" for the
# in-flight items from "In Translation" to "Ready for handoff", refreshes
# the corresponding timestamp columns, and widens the Status columns to
# fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Columns: A File Name | B Path And Name | C Extension | D Publish URL |
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 01:00:21"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn detail sheet ------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 01:00:08"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de detail sheet ------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 01:00:21"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
